# The deck currently applies the "Integral" design (the colour scheme in
# use - ppt/theme/theme2.xml, the theme referenced by the slide master and
# the presentation itself - holds the green/teal "Integral" palette). The
# edit switches the active colour scheme back to the built-in "Office"
# palette, i.e. exactly what you get from Design > Variants > Colors >
# "Office" in the ribbon.
#
# Every slide shares the same ThemeColorScheme object (they all follow the
# single Slide Master), so editing it once via slide 1 updates the shared
# theme part for the whole deck.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# Office theme colours, in ThemeColorScheme order (value = BGR-packed long,
# i.e. the same encoding PowerPoint's ColorFormat.RGB always uses):
#  1 Dark1, 2 Light1, 3 Dark2, 4 Light2,
#  5-10 Accent1-6, 11 Hyperlink, 12 FollowedHyperlink
$tcs.Item(1).RGB  = 0         # dk1      000000
$tcs.Item(2).RGB  = 16777215  # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388   # dk2      44546A
$tcs.Item(4).RGB  = 15132391  # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939  # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501   # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845  # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407     # accent4  FFC000
$tcs.Item(9).RGB  = 12874308  # accent5  4472C4
$tcs.Item(10).RGB = 4697456   # accent6  70AD47
$tcs.Item(11).RGB = 12673797  # hlink    0563C1
$tcs.Item(12).RGB = 7491477   # folHlink 954F72
